# [Fonds de solidarite] Add 2020-12-16 data
#
# The source sheet stores nombre_aides / nombre_entreprises / montant_total
# as TEXT (inline strings), not numbers. A plain Range.Value assignment of a
# numeric-looking string gets auto-coerced to a Number by Excel, which would
# change the cell type. Prefixing the literal with a leading apostrophe
# forces Excel to keep it as text (matching the original t="inlineStr"
# representation) without the apostrophe itself ending up in the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bourgogne-Franche-Comte (reg 27)
$ws.Range("C10").Value = "'496"
$ws.Range("D10").Value = "'441"
$ws.Range("E10").Value = "'2890139.57"

$ws.Range("C11").Value = "'227"
$ws.Range("E11").Value = "'2008003.64"

$ws.Range("C12").Value = "'70"
$ws.Range("E12").Value = "'1021337.89"

# Grand Est (reg 44)
$ws.Range("C36").Value = "'794"
$ws.Range("E36").Value = "'3815114.93"

$ws.Range("C37").Value = "'385"
$ws.Range("E37").Value = "'3324428.91"

$ws.Range("C38").Value = "'150"
$ws.Range("E38").Value = "'2102100.47"

$ws.Range("C39").Value = "'61"
$ws.Range("E39").Value = "'1221474.34"

$ws.Range("C40").Value = "'8"
$ws.Range("E40").Value = "'189000.00"

# Hauts-de-France (reg 32)
$ws.Range("C56").Value = "'1026"
$ws.Range("E56").Value = "'5835085.93"

$ws.Range("C57").Value = "'512"
$ws.Range("E57").Value = "'4821219.81"

$ws.Range("C58").Value = "'189"
$ws.Range("E58").Value = "'2070594.11"

$ws.Range("C59").Value = "'68"
$ws.Range("E59").Value = "'1055181.06"

# Ile-de-France (reg 11)
$ws.Range("C63").Value = "'5714"
$ws.Range("E63").Value = "'24444042.70"

$ws.Range("C64").Value = "'3147"
$ws.Range("E64").Value = "'19420486.38"

$ws.Range("C67").Value = "'50"
$ws.Range("E67").Value = "'1842881.18"

# La Reunion (reg 04)
$ws.Range("C74").Value = "'4"
$ws.Range("D74").Value = "'4"
$ws.Range("E74").Value = "'140000.00"
